# "Refining number labels and the masterlist."
#
# The sheet holds a 1..8-column "masterlist" of sequential numbers (A-H),
# each column continuing straight on from the last row of the column
# before it. This extends that same pattern two columns further, into I
# and J, and nudges the saved selection/column width to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give the new columns the same look as column H (bold/underlined
# header font, thin box border, centered) by copying its formatting. ---
$ws.Range("H1:H45").Copy()
$ws.Range("I1:J45").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column I continues straight on from the bottom of column H. ---
$ws.Range("I1").Formula = "=H45+1"
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=I" + ($r - 1) + "+1"
}

# --- Column J continues straight on from the bottom of column I. ---
$ws.Range("J1").Formula = "=I45+1"
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 10).Formula = "=J" + ($r - 1) + "+1"
}

# --- Column I picked up a slightly narrower custom width. ---
$ws.Range("I1").ColumnWidth = 9.83

# --- The author's last selection before saving. ---
$ws.Range("F7").Select() | Out-Null
